$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string for roc_auc column (Z2:Z17) - all cells must be set
# to the same new string so it replaces the single shared-string entry.
$rocAuc = "[0.66985646 0.64114833 0.54066986 0.60287081 0.67      ]"
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("Z$r").Value = $rocAuc
}

# Row 2
$ws.Range("B2").Value = [double]"0.004255390167236328"
$ws.Range("C2").Value = [double]"0.002134899894764335"
$ws.Range("D2").Value = [double]"0.001120805740356445"
$ws.Range("E2").Value = [double]"0.0002697294199972461"
$ws.Range("R2").Value = [double]"0.7"
$ws.Range("S2").Value = [double]"0.7500000000000001"
$ws.Range("T2").Value = [double]"0.5405405405405405"
$ws.Range("U2").Value = [double]"0.6666666666666667"
$ws.Range("V2").Value = [double]"0.6666666666666665"
$ws.Range("W2").Value = [double]"0.6647747747747748"
$ws.Range("X2").Value = [double]"0.06920326433411365"
$ws.Range("Y2").Value = 6
$ws.Range("AA2").Value = [double]"0.6249090909090909"
$ws.Range("AB2").Value = [double]"0.04878340505235926"

# Row 3
$ws.Range("B3").Value = [double]"0.003809595108032227"
$ws.Range("C3").Value = [double]"0.001804851336569589"
$ws.Range("D3").Value = [double]"0.0009220600128173828"
$ws.Range("E3").Value = [double]"0.0001166264733287864"
$ws.Range("R3").Value = [double]"0.717948717948718"
$ws.Range("S3").Value = [double]"0.631578947368421"
$ws.Range("T3").Value = [double]"0.5142857142857142"
$ws.Range("U3").Value = [double]"0.5161290322580646"
$ws.Range("V3").Value = [double]"0.5294117647058824"
$ws.Range("W3").Value = [double]"0.5818708353133599"
$ws.Range("X3").Value = [double]"0.08078338181691065"
$ws.Range("Y3").Value = 12
$ws.Range("AA3").Value = [double]"0.6249090909090909"
$ws.Range("AB3").Value = [double]"0.04878340505235926"

# Row 4
$ws.Range("B4").Value = [double]"0.001779794692993164"
$ws.Range("C4").Value = [double]"0.001009151129071246"
$ws.Range("D4").Value = [double]"0.0008942127227783203"
$ws.Range("E4").Value = [double]"0.0004473625533180876"
$ws.Range("R4").Value = [double]"0.761904761904762"
$ws.Range("S4").Value = [double]"0.7804878048780488"
$ws.Range("T4").Value = [double]"0.631578947368421"
$ws.Range("U4").Value = [double]"0.7027027027027027"
$ws.Range("V4").Value = [double]"0.7692307692307692"
$ws.Range("W4").Value = [double]"0.7291809972169407"
$ws.Range("X4").Value = [double]"0.05573953308217269"
$ws.Range("Y4").Value = 1
$ws.Range("AA4").Value = [double]"0.6249090909090909"
$ws.Range("AB4").Value = [double]"0.04878340505235926"

# Row 5
$ws.Range("B5").Value = [double]"0.00194544792175293"
$ws.Range("C5").Value = [double]"0.0005363525825055996"
$ws.Range("D5").Value = [double]"0.001091480255126953"
$ws.Range("E5").Value = [double]"0.0007875697058363001"
$ws.Range("R5").Value = [double]"0.7317073170731707"
$ws.Range("S5").Value = [double]"0.6666666666666667"
$ws.Range("T5").Value = [double]"0.6153846153846154"
$ws.Range("U5").Value = [double]"0.7222222222222222"
$ws.Range("V5").Value = [double]"0.5555555555555556"
$ws.Range("W5").Value = [double]"0.6583072753804462"
$ws.Range("X5").Value = [double]"0.06624626274899952"
$ws.Range("Y5").Value = 7
$ws.Range("AA5").Value = [double]"0.6249090909090909"
$ws.Range("AB5").Value = [double]"0.04878340505235926"

# Row 6
$ws.Range("B6").Value = [double]"0.002114534378051758"
$ws.Range("C6").Value = [double]"0.001056559135856211"
$ws.Range("D6").Value = [double]"0.0008840560913085938"
$ws.Range("E6").Value = [double]"0.0003701165360454291"
$ws.Range("R6").Value = [double]"0.7142857142857143"
$ws.Range("S6").Value = [double]"0.717948717948718"
$ws.Range("T6").Value = [double]"0.631578947368421"
$ws.Range("U6").Value = [double]"0.6666666666666667"
$ws.Range("V6").Value = [double]"0.6842105263157895"
$ws.Range("W6").Value = [double]"0.6829381145170619"
$ws.Range("X6").Value = [double]"0.03197673360754662"
$ws.Range("Y6").Value = 4
$ws.Range("AA6").Value = [double]"0.6249090909090909"
$ws.Range("AB6").Value = [double]"0.04878340505235926"

# Row 7
$ws.Range("B7").Value = [double]"0.00311436653137207"
$ws.Range("C7").Value = [double]"0.001747063969683034"
$ws.Range("D7").Value = [double]"0.0007877349853515625"
$ws.Range("E7").Value = [double]"0.0001448485779627544"
$ws.Range("R7").Value = [double]"0.717948717948718"
$ws.Range("S7").Value = [double]"0.631578947368421"
$ws.Range("T7").Value = [double]"0.5142857142857142"
$ws.Range("U7").Value = [double]"0.5161290322580646"
$ws.Range("V7").Value = [double]"0.5294117647058824"
$ws.Range("W7").Value = [double]"0.5818708353133599"
$ws.Range("X7").Value = [double]"0.08078338181691065"
$ws.Range("Y7").Value = 12
$ws.Range("AA7").Value = [double]"0.6249090909090909"
$ws.Range("AB7").Value = [double]"0.04878340505235926"

# Row 8
$ws.Range("B8").Value = [double]"0.00181884765625"
$ws.Range("C8").Value = [double]"0.0008195693556861795"
$ws.Range("D8").Value = [double]"0.0007009506225585938"
$ws.Range("E8").Value = [double]"7.711593565371269E-05"
$ws.Range("R8").Value = [double]"0.7142857142857143"
$ws.Range("S8").Value = [double]"0.761904761904762"
$ws.Range("T8").Value = [double]"0.6500000000000001"
$ws.Range("U8").Value = [double]"0.7027027027027027"
$ws.Range("V8").Value = [double]"0.717948717948718"
$ws.Range("W8").Value = [double]"0.7093683793683795"
$ws.Range("X8").Value = [double]"0.03585222335390306"
$ws.Range("Y8").Value = 2
$ws.Range("AA8").Value = [double]"0.6249090909090909"
$ws.Range("AB8").Value = [double]"0.04878340505235926"

# Row 9
$ws.Range("B9").Value = [double]"0.001537609100341797"
$ws.Range("C9").Value = [double]"0.0001937005859222918"
$ws.Range("D9").Value = [double]"0.0006522655487060547"
$ws.Range("E9").Value = [double]"3.259013666245423E-05"
$ws.Range("R9").Value = [double]"0.7317073170731707"
$ws.Range("S9").Value = [double]"0.6666666666666667"
$ws.Range("T9").Value = [double]"0.6153846153846154"
$ws.Range("U9").Value = [double]"0.7222222222222222"
$ws.Range("V9").Value = [double]"0.5555555555555556"
$ws.Range("W9").Value = [double]"0.6583072753804462"
$ws.Range("X9").Value = [double]"0.06624626274899952"
$ws.Range("Y9").Value = 7
$ws.Range("AA9").Value = [double]"0.6249090909090909"
$ws.Range("AB9").Value = [double]"0.04878340505235926"

# Row 10
$ws.Range("B10").Value = [double]"0.001807165145874023"
$ws.Range("C10").Value = [double]"0.0001513579388838597"
$ws.Range("D10").Value = [double]"0.0007008075714111328"
$ws.Range("E10").Value = [double]"4.664057321818651E-05"
$ws.Range("R10").Value = [double]"0.7567567567567567"
$ws.Range("S10").Value = [double]"0.631578947368421"
$ws.Range("T10").Value = [double]"0.5555555555555555"
$ws.Range("U10").Value = [double]"0.5806451612903226"
$ws.Range("V10").Value = [double]"0.5142857142857143"
$ws.Range("W10").Value = [double]"0.607764427051354"
$ws.Range("X10").Value = [double]"0.08362456130537084"
$ws.Range("Y10").Value = 11
$ws.Range("AA10").Value = [double]"0.6249090909090909"
$ws.Range("AB10").Value = [double]"0.04878340505235926"

# Row 11
$ws.Range("B11").Value = [double]"0.001755952835083008"
$ws.Range("C11").Value = [double]"0.0002532565599078225"
$ws.Range("D11").Value = [double]"0.000637674331665039"
$ws.Range("E11").Value = [double]"2.436381249994812E-05"
$ws.Range("R11").Value = [double]"0.717948717948718"
$ws.Range("S11").Value = [double]"0.631578947368421"
$ws.Range("T11").Value = [double]"0.5142857142857142"
$ws.Range("U11").Value = [double]"0.5161290322580646"
$ws.Range("V11").Value = [double]"0.5294117647058824"
$ws.Range("W11").Value = [double]"0.5818708353133599"
$ws.Range("X11").Value = [double]"0.08078338181691065"
$ws.Range("Y11").Value = 12
$ws.Range("AA11").Value = [double]"0.6249090909090909"
$ws.Range("AB11").Value = [double]"0.04878340505235926"

# Row 12
$ws.Range("B12").Value = [double]"0.001327180862426758"
$ws.Range("C12").Value = [double]"0.0001795394547505557"
$ws.Range("D12").Value = [double]"0.0006309032440185546"
$ws.Range("E12").Value = [double]"3.041548340657028E-05"
$ws.Range("R12").Value = [double]"0.7804878048780488"
$ws.Range("S12").Value = [double]"0.7"
$ws.Range("T12").Value = [double]"0.5789473684210527"
$ws.Range("U12").Value = [double]"0.742857142857143"
$ws.Range("V12").Value = [double]"0.7"
$ws.Range("W12").Value = [double]"0.7004584632312488"
$ws.Range("X12").Value = [double]"0.06777567424597226"
$ws.Range("Y12").Value = 3
$ws.Range("AA12").Value = [double]"0.6249090909090909"
$ws.Range("AB12").Value = [double]"0.04878340505235926"

# Row 13
$ws.Range("B13").Value = [double]"0.001394176483154297"
$ws.Range("C13").Value = [double]"0.0001151481866534802"
$ws.Range("D13").Value = [double]"0.0006142616271972656"
$ws.Range("E13").Value = [double]"1.116919158975634E-05"
$ws.Range("R13").Value = [double]"0.7317073170731707"
$ws.Range("S13").Value = [double]"0.6666666666666667"
$ws.Range("T13").Value = [double]"0.6153846153846154"
$ws.Range("U13").Value = [double]"0.7222222222222222"
$ws.Range("V13").Value = [double]"0.5555555555555556"
$ws.Range("W13").Value = [double]"0.6583072753804462"
$ws.Range("X13").Value = [double]"0.06624626274899952"
$ws.Range("Y13").Value = 7
$ws.Range("AA13").Value = [double]"0.6249090909090909"
$ws.Range("AB13").Value = [double]"0.04878340505235926"

# Row 14
$ws.Range("B14").Value = [double]"0.002327537536621094"
$ws.Range("C14").Value = [double]"0.001164663583200146"
$ws.Range("D14").Value = [double]"0.0009757041931152344"
$ws.Range("E14").Value = [double]"0.0007048918993664734"
$ws.Range("R14").Value = [double]"0.717948717948718"
$ws.Range("S14").Value = [double]"0.631578947368421"
$ws.Range("T14").Value = [double]"0.5142857142857142"
$ws.Range("U14").Value = [double]"0.5161290322580646"
$ws.Range("V14").Value = [double]"0.5294117647058824"
$ws.Range("W14").Value = [double]"0.5818708353133599"
$ws.Range("X14").Value = [double]"0.08078338181691065"
$ws.Range("Y14").Value = 12
$ws.Range("AA14").Value = [double]"0.6249090909090909"
$ws.Range("AB14").Value = [double]"0.04878340505235926"

# Row 15
$ws.Range("B15").Value = [double]"0.002364253997802734"
$ws.Range("C15").Value = [double]"0.001008450657391694"
$ws.Range("D15").Value = [double]"0.0006611347198486328"
$ws.Range("E15").Value = [double]"8.830568433664599E-05"
$ws.Range("R15").Value = [double]"0.717948717948718"
$ws.Range("S15").Value = [double]"0.631578947368421"
$ws.Range("T15").Value = [double]"0.5142857142857142"
$ws.Range("U15").Value = [double]"0.5161290322580646"
$ws.Range("V15").Value = [double]"0.5294117647058824"
$ws.Range("W15").Value = [double]"0.5818708353133599"
$ws.Range("X15").Value = [double]"0.08078338181691065"
$ws.Range("Y15").Value = 12
$ws.Range("AA15").Value = [double]"0.6249090909090909"
$ws.Range("AB15").Value = [double]"0.04878340505235926"

# Row 16
$ws.Range("B16").Value = [double]"0.001373100280761719"
$ws.Range("C16").Value = [double]"0.0001090106895864201"
$ws.Range("D16").Value = [double]"0.0006118297576904296"
$ws.Range("E16").Value = [double]"1.415523852669371E-05"
$ws.Range("R16").Value = [double]"0.7500000000000001"
$ws.Range("S16").Value = [double]"0.6666666666666667"
$ws.Range("T16").Value = [double]"0.631578947368421"
$ws.Range("U16").Value = [double]"0.7222222222222222"
$ws.Range("V16").Value = [double]"0.5555555555555556"
$ws.Range("W16").Value = [double]"0.6652046783625731"
$ws.Range("X16").Value = [double]"0.06869789546474113"
$ws.Range("Y16").Value = 5
$ws.Range("AA16").Value = [double]"0.6249090909090909"
$ws.Range("AB16").Value = [double]"0.04878340505235926"

# Row 17
$ws.Range("B17").Value = [double]"0.001378250122070313"
$ws.Range("C17").Value = [double]"0.0001031700196930825"
$ws.Range("D17").Value = [double]"0.0006153583526611328"
$ws.Range("E17").Value = [double]"1.802734573671376E-05"
$ws.Range("R17").Value = [double]"0.7317073170731707"
$ws.Range("S17").Value = [double]"0.6666666666666667"
$ws.Range("T17").Value = [double]"0.6153846153846154"
$ws.Range("U17").Value = [double]"0.7222222222222222"
$ws.Range("V17").Value = [double]"0.5555555555555556"
$ws.Range("W17").Value = [double]"0.6583072753804462"
$ws.Range("X17").Value = [double]"0.06624626274899952"
$ws.Range("Y17").Value = 7
$ws.Range("AA17").Value = [double]"0.6249090909090909"
$ws.Range("AB17").Value = [double]"0.04878340505235926"
